# Insert a new data row at row 185 (pushing the existing rows 185:233 down
# to 186:234) and populate it with the new record, matching the fixed
# (constant-across-dataset) columns used by every other row in this table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 185:233 down by one row.
$ws.Rows("185:185").Insert()

# Fill in the new row 185 with its data.
$ws.Cells.Item(185, 1).Value = 3
$ws.Cells.Item(185, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(185, 3).Value = "Coquimbo"
$ws.Cells.Item(185, 4).Value = 44855
$ws.Cells.Item(185, 5).Value = 5
$ws.Cells.Item(185, 6).Value = 100112010
$ws.Cells.Item(185, 7).Value = "Achicoria"
$ws.Cells.Item(185, 8).Value = "Sin especificar"
$ws.Cells.Item(185, 9).Value = "Primera"
$ws.Cells.Item(185, 10).Value = 65
$ws.Cells.Item(185, 11).Value = 5000
$ws.Cells.Item(185, 12).Value = 5000
$ws.Cells.Item(185, 13).Value = 5000
$ws.Cells.Item(185, 14).Value = '$/caja 16 unidades'
$ws.Cells.Item(185, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(185, 16).Value = 312
$ws.Cells.Item(185, 17).Value = 16
$ws.Cells.Item(185, 18).Value = "Hortaliza"
